# Restructure the "While ( / If ( i > N ) / Print "Max"" block:
#  - "While  (" keeps its own paragraph but loses the _GoBack bookmark
#  - "If ( i > N )" paragraph is unchanged
#  - "Print "Max"" paragraph becomes two runs ("Print " then "Max", curly
#    quotes removed) with the _GoBack bookmark relocated between them
$d = $word.ActiveDocument

# Locate the "While  (" paragraph robustly via Find rather than a
# hard-coded paragraph index.
$findRange = $d.Content
$found = $findRange.Find.Execute("While  (")
if (-not $found) {
    throw "Could not find 'While  (' paragraph"
}

$allParas = $d.Paragraphs
$whileIdx = -1
for ($i = 1; $i -le $allParas.Count; $i++) {
    $pr = $allParas.Item($i).Range
    if ($findRange.Start -ge $pr.Start -and $findRange.Start -lt $pr.End) {
        $whileIdx = $i
        break
    }
}
if ($whileIdx -eq -1) {
    throw "Could not resolve paragraph index for 'While  (' "
}

$paraWhile = $allParas.Item($whileIdx)
$paraIf = $allParas.Item($whileIdx + 1)
$paraPrint = $allParas.Item($whileIdx + 2)

# Sanity-check the paragraphs we are about to replace.
if ($paraIf.Range.Text -notmatch "If \( i > N \)") {
    throw "Unexpected paragraph after 'While  (': $($paraIf.Range.Text)"
}
if ($paraPrint.Range.Text -notmatch "Print") {
    throw "Unexpected paragraph after 'If ( i > N )': $($paraPrint.Range.Text)"
}

$target = $d.Range($paraWhile.Range.Start, $paraPrint.Range.End)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00761337" w:rsidRPr="00FF6E6F" w:rsidRDefault="00761337" w:rsidP="00201761"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r w:rsidRPr="00FF6E6F"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>While  (</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:tab/><w:t>If ( i &gt; N )</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:tab/><w:t xml:space="preserve">Print </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Max</w:t></w:r></w:p>'

$target.InsertXML($xml)
